$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (the "count" column) rows 2-21 change from 1154 to 1279
$ws.Range("A2:A21").Value = 1279
